# Apply cell updates per the commit diff (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.633.19"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3
$ws.Range("D3").Value = "3.675.55"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'649.03"
$ws.Range("E5").Value = "  -4.32%  "

# Row 6
$ws.Range("D6").Value = "'160.72"
$ws.Range("E6").Value = "  -1.19%  "

# Row 7
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  -2.66%  "

# Row 10
$ws.Range("D10").Value = "'7.15"
$ws.Range("E10").Value = "  -0.11%  "

# Row 11
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").Value = "'0.0000230"
$ws.Range("E12").Value = "  -2.45%  "

# Row 13
$ws.Range("D13").Value = "4.295.46"
$ws.Range("E13").Value = "  -0.87%  "

# Row 14
$ws.Range("D14").Value = "'32.60"
$ws.Range("E14").Value = "  -1.11%  "

# Row 15
$ws.Range("D15").Value = "3.681.32"
$ws.Range("E15").Value = "  -0.89%  "

# Row 16
$ws.Range("D16").Value = "69.645.25"
$ws.Range("E16").Value = "  -0.20%  "

# Row 17
$ws.Range("E17").Value = "  +0.44%  "

# Row 18
$ws.Range("D18").Value = "'6.52"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("D19").Value = "'15.88"
$ws.Range("E19").Value = "  -1.66%  "

# Row 20
$ws.Range("E20").Value = "  +4.31%  "

# Row 21
$ws.Range("D21").Value = "'469.43"
$ws.Range("E21").Value = "  -0.96%  "

# Row 22
$ws.Range("D22").Value = "'0.655"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("E23").Value = "  -1.14%  "

# Row 24
$ws.Range("D24").Value = "3.824.41"
$ws.Range("E24").Value = "  -0.78%  "

# Row 25
$ws.Range("E25").Value = "  -0.11%  "

# Row 26
$ws.Range("D26").Value = "'0.0000125"
$ws.Range("E26").Value = "  -2.95%  "

# Row 27
$ws.Range("D27").Value = "'11.13"
$ws.Range("E27").Value = "  +0.87%  "

# Row 28
$ws.Range("D28").Value = "'8.78"
$ws.Range("E28").Value = "  -4.28%  "

# Row 29
$ws.Range("D29").Value = "'2.64"
$ws.Range("E29").Value = "  -2.72%  "

# Row 30
$ws.Range("E30").Value = "  -3.78%  "

# Row 31
$ws.Range("E31").Value = "  -0.15%  "

# Row 32
$ws.Range("E32").Value = "  -2.35%  "

# Row 33
$ws.Range("E33").Value = "  -2.18%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.70"
$ws.Range("E34").Value = "  -1.19%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.165"
$ws.Range("E35").Value = "  -0.79%  "

# Row 36
$ws.Range("D36").Value = "3.670.63"
$ws.Range("E36").Value = "  -0.70%  "

# Row 37
$ws.Range("D37").Value = "'8.36"
$ws.Range("E37").Value = "  -2.73%  "

# Row 39
$ws.Range("E39").Value = "  -5.65%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'177.70"
$ws.Range("E40").Value = "  +6.11%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.09%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -2.41%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0895"
$ws.Range("E43").Value = "  -1.57%  "

# Row 44
$ws.Range("E44").Value = "  -1.97%  "

# Row 45
$ws.Range("D45").Value = "'46.70"
$ws.Range("E45").Value = "  -0.65%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'28.97"
$ws.Range("E46").Value = "  +2.51%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.78"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48
$ws.Range("D48").Value = "'0.000268"
$ws.Range("E48").Value = "  -5.28%  "

# Row 49
$ws.Range("D49").Value = "'7.83"
$ws.Range("E49").Value = "  -1.32%  "

# Row 50
$ws.Range("E50").Value = "  -4.97%  "

# Row 51
$ws.Range("D51").Value = "'1.04"
$ws.Range("E51").Value = "  -6.51%  "
